# Add new column 'Crrection' (column P) to the Card1 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card1")

# --- Header row: P1 gets the same (bold/border/centered) header style as
# the rest of row 1, copied from the neighboring header cell O1 so we
# reuse the existing style definition instead of creating a new one.
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(1, 16).Value = "Crrection"

# --- Column M ("Event") rows 2-12 were blank; they now carry the
# placeholder text "nan" like the other data columns.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"
}

# --- New column P rows 2-12: create the (empty) cells so the column
# exists all the way down, matching the unstyled data cells next to them.
for ($r = 2; $r -le 12; $r++) {
    $ws.Range("O$r").Copy()
    $ws.Range("P$r").PasteSpecial(-4122)   # xlPasteFormats
}

$excel.CutCopyMode = 0
